$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "MCF"
$ws.Range("B40").Value = "Modification Consent Form"
$ws.Range("C40").Value = "Modification Consent Form"
$ws.Range("D40").Value = "eng"
$ws.Range("E40").Value = "t"
$ws.Range("F40").Value = "zimbe"

$ws.Range("G39").Copy()
$ws.Range("G40").PasteSpecial(-4122)
$ws.Range("G40").Value = 45634.747916666667

$ws.Rows("40:40").Select()
